$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) figures.
# NumberFormat "@" + resetting Style to "Normal" afterwards keeps the
# cell as plain text (matching the source inlineStr cells) instead of
# letting Excel auto-convert numeric-looking text into a Number/Double.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '44.010.01'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.65%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.240.28'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -2.01%  '
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.26%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '314.79'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -1.67%  '
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '98.81'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -8.30%  '
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -3.16%  '
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.29%  '
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.532'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -7.69%  '
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '36.09'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -8.17%  '
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0823'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -2.67%  '
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.36'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -7.71%  '
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -2.83%  '
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.580.34'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -1.79%  '
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.841'
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.241.45'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -1.52%  '
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '13.91'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -5.58%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '43.820.69'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -0.76%  '
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.21'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -6.72%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.0₃0973'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -3.27%  '
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.31'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -4.26%  '
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '65.92'
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '237.73'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -0.64%  '
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -8.13%  '
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.02'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -9.26%  '
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +0.25%  '
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.12'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -1.72%  '
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -4.26%  '
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '36.55'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -7.09%  '
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '5.98'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -9.15%  '
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '20.04'
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -3.42%  '
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '154.20'
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -5.68%  '
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0833'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -6.76%  '
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.34'
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +1.12%  '
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -2.84%  '
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.92'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -7.96%  '
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.108'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -7.00%  '
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -3.34%  '
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '15.62'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +1.41%  '
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '4.00'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -11.32%  '
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -13.03%  '
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0307'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -7.02%  '
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.707.25'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -3.84%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '82.46'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -4.58%  '
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -6.75%  '
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -5.78%  '
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '101.63'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -3.32%  '
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '71.50'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -5.80%  '
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '56.20'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -6.79%  '
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.61'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -6.24%  '
$c.Style = "Normal"
